# BAU Global Solar and Wind Cap.xlsx - EPS v1.5.0 update
# Adds three new generation-technology rows (crude oil, heavy or residual
# fuel oil, municipal solid waste) to the BGSaWC sheet, each initialized to
# zero capacity across all forecast years, plus a header label for the
# existing "Global Capacity (MW)" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BGSaWC")

# New rows 15-17: technology label in column A, zero values for every
# forecast year column B:AI (2017-2050).
$ws.Range("A15").Value = "crude oil"
$ws.Range("B15:AI15").Value = 0

$ws.Range("A16").Value = "heavy or residual fuel oil"
$ws.Range("B16:AI16").Value = 0

$ws.Range("A17").Value = "municipal solid waste"
$ws.Range("B17:AI17").Value = 0

# Label the year header row (A1) that previously had no caption.
$ws.Range("A1").Value = "Global Capacity (MW)"

# Restore the selection to A2, matching the saved workbook state.
$ws.Range("A2").Select() | Out-Null
